$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Resize the two columns of the first table (address / date block).
#    6237/3686 dxa (twips) -> 5529/4394 dxa == 276.45/219.7 points.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Width = 276.45
$t.Cell(1, 2).Width = 219.7

# ---------------------------------------------------------------------------
# 2) Split the "u.p Direktur" run into proofed sub-runs, matching what Word's
#    spelling/grammar checker inserts around "u.p" and "Direktur".
# ---------------------------------------------------------------------------
$find = $d.Content
$found = $find.Find.Execute("u.p Direktur", $true, $false, $false, $false, $false, `
                             $true, 1, $false, "", 0)

if ($found) {
    # Use a brand-new Range object addressing the same span: reusing the
    # Find range after further manipulation can leave it referencing stale
    # match state.
    $target = $d.Range($find.Start, $find.End)

    $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p>' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:proofErr w:type="gramStart"/>' +
           '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t>u.p</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '<w:proofErr w:type="gramEnd"/>' +
           '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
           '<w:proofErr w:type="spellStart"/>' +
           '<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t>Direktur</w:t></w:r>' +
           '<w:proofErr w:type="spellEnd"/>' +
           '</w:p>' +
           '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData>' +
           '</pkg:part>' +
           '</pkg:package>'

    $target.InsertXML($xml)
}
